# Complete restructure and rewrite of documentation ready for v2
$wb = $excel.ActiveWorkbook

$notes   = $wb.Worksheets.Item("Notes")
$studies = $wb.Worksheets.Item("studies")
$surveys = $wb.Worksheets.Item("surveys")
$counts  = $wb.Worksheets.Item("counts")

# ---------------------------------------------------------------------------
# Sheet "studies" - new header schema, sample row updated
# ---------------------------------------------------------------------------
$studies.Range("A1").Value = "study_id"
$studies.Range("B1").Value = "study_label"
$studies.Range("C1").Value = "description"
$studies.Range("D1").Value = "access_level"
$studies.Range("E1").Value = "contributors"
$studies.Range("F1").Value = "reference"
$studies.Range("G1").Value = "reference_year"

$studies.Range("A2").Value = "foo"
$studies.Range("C2").ClearContents()
$studies.Range("D2").Value = "public"
# F2 keeps its existing hyperlinked URL value/style untouched

# ---------------------------------------------------------------------------
# Sheet "surveys" - two new columns (location_method/location_notes and
# time_method) inserted into the middle of the header row, sample row values
# re-homed to match
# ---------------------------------------------------------------------------

# Capture the row-2 values that need to shift into new column positions
# before the columns move around.
$oldG2 = $surveys.Range("G2").Value()
$oldJ2 = $surveys.Range("J2").Value()
$oldK2 = $surveys.Range("K2").Value()

# Vacate the old positions completely (value + any baked-in formatting)
$surveys.Range("G2").Clear()
$surveys.Range("J2").Clear()
$surveys.Range("K2").Clear()

$surveys.Range("A2").Value = "foo"

# Re-home the moved sample values
$surveys.Range("H2").Value = $oldG2
$surveys.Range("H2").Style = "Normal"

$surveys.Range("K2").NumberFormat = "@"
$surveys.Range("K2").Value = $oldJ2

$surveys.Range("M2").Value = $oldK2
$surveys.Range("M2").Style = "Normal"

# New column with no sample data, but still carrying the "Text" format
$surveys.Range("L2").NumberFormat = "@"

# New header row
$surveys.Range("A1").Value = "study_id"
$surveys.Range("B1").Value = "survey_id"
$surveys.Range("C1").Value = "country_name"
$surveys.Range("D1").Value = "site_name"
$surveys.Range("E1").Value = "latitude"
$surveys.Range("F1").Value = "longitude"
$surveys.Range("G1").Value = "location_method"
$surveys.Range("H1").Value = "location_notes"
$surveys.Range("I1").Value = "collection_start"
$surveys.Range("J1").Value = "collection_end"
$surveys.Range("K1").Value = "collection_day"
$surveys.Range("L1").Value = "time_method"
$surveys.Range("M1").Value = "time_notes"

# Headers drop the old special header font entirely; I1:L1 keep the "Text"
# number format that was already used for the date-like columns
$surveys.Range("A1:H1").Style = "Normal"
$surveys.Range("M1").Style = "Normal"
$surveys.Range("I1:L1").Style = "Normal"
$surveys.Range("I1:L1").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Sheet "counts" - header renamed, sample key renamed (study01 -> foo)
# ---------------------------------------------------------------------------
$counts.Range("A1").Value = "study_id"
$counts.Range("B1").Value = "survey_id"
$counts.Range("C1").Value = "variant_string"
$counts.Range("D1").Value = "variant_num"
$counts.Range("E1").Value = "total_num"

$counts.Range("A2").Value = "foo"
$counts.Range("A3").Value = "foo"

# ---------------------------------------------------------------------------
# View state: "studies" becomes the active/selected tab, and each sheet's
# cursor selection moves to reflect where editing left off
# ---------------------------------------------------------------------------
$studies.Activate()
$studies.Range("D3").Select()

$surveys.Activate()
$surveys.Range("A1:M2").Select()

$counts.Activate()
$counts.Range("D8").Select()

$notes.Activate()
$notes.Range("A3").Select()

$studies.Activate()
